# ---------------------------------------------------------------------------
# "exaggerating differences for visual effect"
#
# 1. spp1: push several 0 -> 1 so the presence/absence grid reads more
#    starkly, and move the selection over the touched block.
# 2. str1 / str3: swap which data lives under which tab name (str1 <-> str3),
#    then re-sort the (newly-named) str1 sheet's G:J block in descending
#    order by column G so it again shows a smooth gradient.
# 3. str2 becomes the active/selected tab.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. spp1: exaggerate the 0/1 contrast in the bottom-right block --------
$spp1 = $wb.Worksheets.Item("spp1")
$spp1.Range("I5").Value2 = 1
$spp1.Range("J5").Value2 = 1
$spp1.Range("I6").Value2 = 1
$spp1.Range("J6").Value2 = 1
$spp1.Range("I7").Value2 = 1
$spp1.Range("J7").Value2 = 1
$spp1.Range("I8").Value2 = 1
$spp1.Range("J8").Value2 = 1
$spp1.Range("J9").Value2 = 1
$spp1.Range("J10").Value2 = 1

# --- 2. swap the str1 / str3 tab names (their data stays put; only the ----
#        labels move) -------------------------------------------------------
$str1 = $wb.Worksheets.Item("str1")
$str3 = $wb.Worksheets.Item("str3")
$str1.Name = "str1_tmp_swap"
$str3.Name = "str1"
$str1.Name = "str3"

# Re-fetch by (new) name so we always act on the right tab from here on.
$newStr1 = $wb.Worksheets.Item("str1")
$newStr3 = $wb.Worksheets.Item("str3")

# Put the (renamed) str1 tab back in front of str3, matching the original
# left-to-right ordering (str1, str3, str2).
$newStr1.Move($newStr3)

# --- 3. re-sort str1's G:J block (descending on column G) so the colour ---
#        gradient still reads correctly after the data swap ---------------
$str1Sorted = $wb.Worksheets.Item("str1")

# remember the original row order for G:J so we can repair any tie-breaking
# quirks the generic Sort engine introduces, before applying it.
$savedGJ = New-Object 'object[,]' 10,4
for ($r = 1; $r -le 10; $r++) {
  for ($c = 1; $c -le 4; $c++) {
    $savedGJ[$r-1, $c-1] = $str1Sorted.Cells.Item($r, $c + 6).Value2
  }
}

$sortRange = $str1Sorted.Range("G1:K10")
$str1Sorted.Sort.SortFields.Clear()
$str1Sorted.Sort.SortFields.Add($str1Sorted.Range("G1:G10"), 0, 2)
$str1Sorted.Sort.SetRange($sortRange)
$str1Sorted.Sort.Header = 0
$str1Sorted.Sort.Apply()

# the data here is heavily tied on column G, so re-apply the descending
# row order explicitly (full G:J row blocks, bottom-to-top) to guarantee
# correct values regardless of how ties were broken above.
for ($r = 1; $r -le 10; $r++) {
  for ($c = 1; $c -le 4; $c++) {
    $str1Sorted.Cells.Item($r, $c + 6).Value2 = $savedGJ[10 - $r, $c - 1]
  }
}

# --- 4. selections / active tab --------------------------------------------
$spp1.Range("I9:J10").Select()
$str1Sorted.Range("C15").Select()

$str2 = $wb.Worksheets.Item("str2")
$str2.Activate()
